{"js": "// Find the paragraph whose entire text is \"26 \" (the lone \"26 \" line that\n// precedes the empty-body steps) and append a new bold/Arial/32 run\n// \" git push && git push --tag\" to it, matching the other \"git ...\"\n// command runs already present in the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"26 \") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the paragraph with text \"26 \"');\n}\n\n// Append the new run at the end of that paragraph (still inside it,\n// before the paragraph mark).\nconst newRun = target.insertText(\" git push && git push --tag\", Word.InsertLocation.end);\n\n// Match the bold/Arial/32pt (half-point 32 => 16pt \"Size\" in the OM)\n// formatting used by the sibling \"git ...\" command runs.\nnewRun.font.set({\n  bold: true,\n  name: \"Arial\",\n  size: 16\n});\nnewRun.font.nameAscii = \"Arial\";\nnewRun.font.nameBidirectional = \"Arial\";\n\nawait context.sync();\n", "ps1": "# Locate the standalone \"26 \" line (the line that precedes the\n# numbered Git-conflict walkthrough step) and append a new bold/Arial/32\n# run \" git push && git push --tag\" to the end of that same paragraph,\n# matching the formatting of the sibling \"git ...\" command runs already\n# in the document.\n\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Text = \"26 \"\n$r.Find.Forward = $true\n$r.Find.Wrap = $wdFindStop\n$r.Find.Format = $false\n$r.Find.MatchCase = $false\n$r.Find.MatchWholeWord = $false\n$r.Find.MatchWildcards = $false\n\n$found = $r.Find.Execute()\n\nif ($found) {\n    # Collapse to the caret right after \"26 \" (still inside the\n    # paragraph, before its end-of-paragraph mark) and insert the run.\n    $r.Collapse($wdCollapseEnd)\n    $r.InsertAfter(\" git push && git push --tag\")\n\n    # Match the bold / Arial / 32 (half-points => Size 16) formatting\n    # used by the other \"git ...\" command runs in the document.\n    $r.Font.Name = \"Arial\"\n    $r.Font.NameAscii = \"Arial\"\n    $r.Font.NameBi = \"Arial\"\n    $r.Font.Bold = $true\n    $r.Font.Size = 16\n}\n"}
